# Apply changes described in the commit:
# "Cambios a la BD en tabla uaprendizaje y p_ua, se elimino la clavessh"
#
# Adds 4 new log rows (69-72) to the Bitacora worksheet describing
# changes made to the uaprendizaje and p_ua tables:
#   69: IVAN | UAPRENDIZAJE | Se elimino el campo caracter... | 05-nov-14 | Se cambio por requerimiento de usuario
#   70: IVAN | UAPRENDIZAJE | Se elimino el campo semestre_sug | 05-nov-14 | Se cambio por requerimiento de usuario
#   71: IVAN | P_UA | Se AGREGO el campo caracter...           | (no date) | Ver registro 69
#   72: IVAN | P_UA | Se AGREGO el campo semestre_sug          | (no date) | Ver registro 70

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# Write the text content of new cells in the order the original
# spreadsheet's shared-string table introduces them, so the resulting
# sharedStrings.xml sequence matches the source file exactly.
$ws.Range("A69").Value = "IVAN"
$ws.Range("B69").Value = "UAPRENDIZAJE"
$ws.Range("C69").Value = "Se elimino el campo carácter y la llave foranea fk_uaprendizaje_caracter1 "

$ws.Range("A70").Value = "IVAN"
$ws.Range("B70").Value = "UAPRENDIZAJE"
$ws.Range("C70").Value = "Se elimino el campo semestre_sug"

$ws.Range("G69").Value = "Se cambio por requerimiento de usuario"
$ws.Range("G70").Value = "Se cambio por requerimiento de usuario"

$ws.Range("A71").Value = "IVAN"
$ws.Range("B71").Value = "P_UA"
$ws.Range("C71").Value = "Se AGREGO el campo carácter y la llave foranea fk_uaprendizaje_caracter1 "

$ws.Range("A72").Value = "IVAN"
$ws.Range("B72").Value = "P_UA"
$ws.Range("C72").Value = "Se AGREGO el campo semestre_sug"

$ws.Range("G71").Value = "Ver registro 69"
$ws.Range("G72").Value = "Ver registro 70"

# Dates for rows 69/70 (05-Nov-2014, serial 41948), formatted like the
# other "FECHA PUBLICACION CAMBIO" cells in column D (built-in format
# d-mmm-yy, centered).
$ws.Range("D69").Value = 41948
$ws.Range("D70").Value = 41948
$ws.Range("D69:D70").NumberFormat = "d-mmm-yy"
$ws.Range("D69:D70").HorizontalAlignment = $xlCenter

# Columns A and B use centered horizontal/vertical alignment throughout
# the log (style used by every other data row).
$ws.Range("A69:B72").HorizontalAlignment = $xlCenter
$ws.Range("A69:B72").VerticalAlignment = $xlCenter

# Leave the active selection on the next empty row, like the author did.
$ws.Range("G73").Select()
